$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.248.98"
$ws.Range("E2").Value = "  +1.93%  "
$ws.Range("D3").Value = "3.192.99"
$ws.Range("E3").Value = "  +1.37%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'596.00"
$ws.Range("E5").Value = "  +3.89%  "
$ws.Range("D6").Value = "'154.20"
$ws.Range("E6").Value = "  +3.08%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "3.194.41"
$ws.Range("E8").Value = "  +1.30%  "
$ws.Range("E10").Value = "  +0.53%  "
$ws.Range("D11").Value = "'6.10"
$ws.Range("E11").Value = "  -0.52%  "
$ws.Range("D12").Value = "'0.513"
$ws.Range("E12").Value = "  +3.32%  "
$ws.Range("D13").Value = "'0.0000270"
$ws.Range("E13").Value = "  +2.82%  "
$ws.Range("D14").Value = "'39.05"
$ws.Range("E14").Value = "  +5.45%  "
$ws.Range("D15").Value = "3.718.38"
$ws.Range("E15").Value = "  +1.45%  "
$ws.Range("D16").Value = "66.157.01"
$ws.Range("E16").Value = "  +1.67%  "
$ws.Range("D17").Value = "'7.43"
$ws.Range("E17").Value = "  +4.75%  "
$ws.Range("D18").Value = "3.199.29"
$ws.Range("E18").Value = "  +1.32%  "
$ws.Range("E19").Value = "  +0.54%  "
$ws.Range("D20").Value = "'509.82"
$ws.Range("E20").Value = "  +0.79%  "
$ws.Range("D21").Value = "'15.29"
$ws.Range("E21").Value = "  +3.62%  "
$ws.Range("D22").Value = "'0.737"
$ws.Range("E22").Value = "  +2.87%  "
$ws.Range("E23").Value = "  -0.47%  "
$ws.Range("D24").Value = "'8.00"
$ws.Range("E24").Value = "  +3.71%  "
$ws.Range("D25").Value = "'84.82"
$ws.Range("E25").Value = "  +0.83%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").Value = "'9.26"
$ws.Range("E27").Value = "  +4.91%  "
$ws.Range("E28").Value = "  +3.21%  "
$ws.Range("D29").Value = "'2.28"
$ws.Range("E29").Value = "  +5.42%  "
$ws.Range("D30").Value = "'6.96"
$ws.Range("E30").Value = "  +12.47%  "
$ws.Range("D31").Value = "'2.90"
$ws.Range("E31").Value = "  +3.77%  "
$ws.Range("D32").Value = "'28.26"
$ws.Range("E32").Value = "  +2.58%  "
$ws.Range("E33").Value = "  +3.41%  "
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("E35").Value = "  +0.59%  "
$ws.Range("D36").Value = "'54.87"
$ws.Range("E36").Value = "  -0.14%  "
$ws.Range("D37").Value = "'0.0900"
$ws.Range("E37").Value = "  +0.50%  "
$ws.Range("D38").Value = "'484.54"
$ws.Range("E38").Value = "  +4.59%  "
$ws.Range("D39").Value = "'0.0417"
$ws.Range("E39").Value = "  -0.59%  "
$ws.Range("E40").Value = "  -2.03%  "
$ws.Range("E41").Value = "  +2.29%  "
$ws.Range("D42").Value = "'0.123"
$ws.Range("E42").Value = "  +5.16%  "
$ws.Range("D43").Value = "'0.297"
$ws.Range("E43").Value = "  +5.77%  "
$ws.Range("D44").Value = "0.0₃0653"
$ws.Range("E44").Value = "  +12.12%  "
$ws.Range("D45").Value = "2.925.95"
$ws.Range("E45").Value = "  -4.06%  "
$ws.Range("D46").Value = "'2.41"
$ws.Range("E46").Value = "  -0.91%  "
$ws.Range("D47").Value = "'28.44"
$ws.Range("E47").Value = "  -0.31%  "
$ws.Range("E48").Value = "  -0.09%  "
$ws.Range("E49").Value = "  +2.19%  "
$ws.Range("D50").Value = "'2.32"
$ws.Range("E50").Value = "  +3.27%  "
$ws.Range("D51").Value = "'2.58"
$ws.Range("E51").Value = "  +6.60%  "
